# Add a new "2021" column (column P) to the table, mirroring the style of
# the existing "2020" column (O) for each row, then fill in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - header separator row, just needs the formatting (no value)
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

# Row 4 - year header
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value2 = 2021

# Row 5 - total, uses style from O8 (s=18)
$ws.Range("O8").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value2 = 9038

# Row 6 - sub-header "By sex:" (empty)
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)

# Row 7 - Women, uses style from O6 (s=17)
$ws.Range("O6").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value2 = 8587

# Row 8 - Men, uses style from O6 (s=17)
$ws.Range("O6").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value2 = 451

# Row 9 - sub-header "By age group:" (empty)
$ws.Range("O6").Copy()
$ws.Range("P9").PasteSpecial(-4122)

# Rows 10-24 - data rows with no 2021 data available, filled with the
# "..." placeholder (same shared string + style as column O)
for ($r = 10; $r -le 24; $r++) {
    $src = $ws.Range("O$r")
    $src.Copy()
    $dst = $ws.Range("P$r")
    $dst.PasteSpecial(-4122)
    $dst.Value2 = $src.Value2
}

# Row 25 - last data row (bottom border), same placeholder pattern
$ws.Range("O25").Copy()
$ws.Range("P25").PasteSpecial(-4122)
$ws.Range("P25").Value2 = $ws.Range("O25").Value2

# Clear clipboard marquee / copy mode
$excel.CutCopyMode = $false

# Match the author's final selection recorded in the saved file
$ws.Range("Q4").Select()
